# Edit: expand the "Messages Metamodel (Context Model):" bullet into a
# full outline with bolded sub-heading runs and explanatory bullet text,
# matching the target OOXML diff.
#
# Strategy: locate the paragraph that holds the title run and the empty
# paragraph that immediately follows it, then replace that whole
# (title-paragraph + blank-paragraph) range in one shot with the exact
# WordprocessingML for: the (reformatted) title paragraph, seven new
# paragraphs (alternating blank / explanatory text, sizes in half-points
# 16 == 8pt), and the trailing blank paragraph (also reformatted).

$d = $word.ActiveDocument

$titleText = "Messages Metamodel (Context Model):"

$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq ($titleText + [char]13)) {
        $titlePara = $cand
        break
    }
}

$nextPara = $titlePara.Next()

$targetRange = $d.Range($titlePara.Range.Start, $nextPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$explainModels = "Explain models (resources, statements, kinds)."
$explainLayers = "Explain layers / aggregation."
$explainMessages = "Explain messages (resource resolution). Grammar. Match model Resource(s). Compound nested CSPO statement contexts defines result behaviors. Message CSPO contexts may define create, retrieve, update or delete operations (passing 'null' for example for resource / statement to be deleted)."
$explainTransforms = "Explain transforms (message appplication). Transform: Resource stream result of Message application over resolved Resource(s)). Input statements: Message(s) / Resource(s) (from input message or to be populated or populated in dialog) and `"goal`" Message / Resource aggregating a model from Resource MetaGraph with Message / Resource bindings."

$xml = '<w:p ' + $ns + '>' +
         '<w:pPr><w:contextualSpacing w:val="0"/>' +
           '<w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="1"/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr>' +
         '</w:pPr>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="1"/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr>' +
           '<w:t xml:space="preserve">' + $titleText + '</w:t></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:b w:val="1"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $explainModels + '</w:t></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $explainLayers + '</w:t></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $explainMessages + '</w:t></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:b w:val="1"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $explainTransforms + '</w:t></w:r>' +
         '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
       '</w:p>' +
       '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr><w:b w:val="1"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
       '</w:p>'

$targetRange.InsertXML($xml)
